$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new report (Congress 103, Session 1)
$ws.Name = "103_1"

# Rewrite the label/value rows for the new confirmations report
$ws.Range("A1").Value = "Labels"
$ws.Range("B1").Value = "Values"
$ws.Range("A2").Value = "Congress"
$ws.Range("B2").Value = 103
$ws.Range("A3").Value = "Session"
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = "Start Date"
$ws.Range("B4").Value = 33974
$ws.Range("A5").Value = "End Date"
$ws.Range("B5").Value = 34297
$ws.Range("A6").Value = "Civilian "
$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("B7").Value = 703
$ws.Range("A8").Value = "     Civilian, Confirmed "
$ws.Range("B8").Value = 510
$ws.Range("A9").Value = "     Civilian, Unconfirmed "
$ws.Range("B9").Value = 172
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("B10").Value = 8
$ws.Range("A11").Value = "     Civilian, Returned at sine die adjournment"
$ws.Range("B11").Value = 13
$ws.Range("A12").Value = "Civilian (FS, PHS, CG, NOAA)"
$ws.Range("A13").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("B13").Value = 2228
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed"
$ws.Range("B14").Value = 2190
$ws.Range("A15").Value = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed "
$ws.Range("B15").Value = 38
$ws.Range("A16").Value = "Air Force "
$ws.Range("A17").Value = "     Air Force, New nominations"
$ws.Range("B17").Value = 13741
$ws.Range("A18").Value = "     Air Force, Confirmed"
$ws.Range("B18").Value = 13736
$ws.Range("A19").Value = "     Air Force, Unconfirmed "
$ws.Range("B19").Value = 4
$ws.Range("A20").Value = "     Air Force, Withdrawn "
$ws.Range("B20").Value = 1
$ws.Range("A21").Value = "Army "
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("B22").Value = 12714
$ws.Range("A23").Value = "     Army, Confirmed"
$ws.Range("B23").Value = 11028
$ws.Range("A24").Value = "     Army, Unconfirmed "
$ws.Range("B24").Value = 1686
$ws.Range("A25").Value = "     Army, Withdrawn "
$ws.Range("B25").Value = 2
$ws.Range("A26").Value = "Navy "
$ws.Range("A27").Value = "     Navy, New nominations"
$ws.Range("B27").Value = 11119
$ws.Range("A28").Value = "     Navy, Confirmed"
$ws.Range("B28").Value = 9584
$ws.Range("A29").Value = "     Navy, Unconfirmed "
$ws.Range("B29").Value = 657
$ws.Range("A30").Value = "     Navy, Withdrawn "
$ws.Range("B30").Value = 878
$ws.Range("A31").Value = "Marine Corps"
$ws.Range("A32").Value = "     Marine Corps, New nominations"
$ws.Range("B32").Value = 1834
$ws.Range("A33").Value = "     Marine Corps, Confirmed "
$ws.Range("B33").Value = 1628
$ws.Range("A34").Value = "     Marine Corps, Unconfirmed "
$ws.Range("B34").Value = 15
$ws.Range("A35").Value = "     Marine Corps, Withdrawn "
$ws.Range("B35").Value = 191
$ws.Range("A36").Value = "Total new nominations"
$ws.Range("B36").Value = 42339
$ws.Range("A37").Value = "Total confirmed "
$ws.Range("B37").Value = 38676
$ws.Range("A38").Value = "Total unconfirmed       "
$ws.Range("B38").Value = 2752
$ws.Range("A39").Value = "Total withdrawn "
$ws.Range("B39").Value = 1078
$ws.Range("A40").Value = "Total returned at sine die adjournment"
$ws.Range("B40").Value = 13

# Match the number formatting used by the other grand-total rows
$ws.Range("B36").NumberFormat = "#,##0"
$ws.Range("B2").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new layout only needs 40 rows; remove the old trailing row 41
$ws.Rows("41:41").Delete()

# Reset the stale selection that pointed past the end of the data
$ws.Range("A1").Select()
